# [JPADCAD] minor changes in sandbox inputs
# Swap the TORENBEEK_1982 / SFORZA rows in the Xcg (and Ycg) estimation
# method comparison tables on the FUSELAGE and WING sheets.

$wb = $excel.ActiveWorkbook

# --- FUSELAGE: Xcg ESTIMATION METHOD COMPARISON (rows 23-24) ---
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("A23").Value = "SFORZA"
$ws.Range("C23").Value = 17.143322222222217
$ws.Range("A24").Value = "TORENBEEK_1982"
$ws.Range("C24").Value = 16.8345

# --- WING: Xcg ESTIMATION METHOD COMPARISON (rows 23-24) ---
$ws = $wb.Worksheets.Item("WING")
$ws.Range("A23").Value = "SFORZA"
$ws.Range("C23").Value = 4.3631082000119275
$ws.Range("A24").Value = "TORENBEEK_1982"
$ws.Range("C24").Value = 3.5939754358446514

# --- WING: Ycg ESTIMATION METHOD COMPARISON (rows 27-28) ---
$ws.Range("A27").Value = "SFORZA"
$ws.Range("C27").Value = 4.998846772296348
$ws.Range("A28").Value = "TORENBEEK_1982"
$ws.Range("C28").Value = 6.114221148470394
